$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain numeric-looking price strings to remain text (matches source data type)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.189.91"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "1.864.75"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "313.11"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "0.5102"
$ws.Range("E7").Value = "  -0.68%  "
$ws.Range("D8").Value = "0.3912"
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("D9").Value = "0.08269"
$ws.Range("E9").Value = "  -1.34%  "
$ws.Range("D10").Value = "1.112"
$ws.Range("E10").Value = "  -0.23%  "
$ws.Range("D11").Value = "6.220"
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("D12").Value = "1.856.84"
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("D13").Value = "20.23"
$ws.Range("E13").Value = "  -1.76%  "
$ws.Range("D14").Value = "7.208"
$ws.Range("E14").Value = "  -1.22%  "
$ws.Range("D15").Value = "1.007"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").Value = "91.28"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").Value = "0.00001097"
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("D18").Value = "0.06683"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").Value = "17.63"
$ws.Range("E19").Value = "  -0.47%  "
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").Value = "5.933"
$ws.Range("E21").Value = "  -1.67%  "
$ws.Range("D22").Value = "28.211.05"
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").Value = "11.04"
$ws.Range("E23").Value = "  -0.46%  "
$ws.Range("D24").Value = "2.253"
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").Value = "2.065.22"
$ws.Range("E25").Value = "  -0.87%  "
$ws.Range("D26").Value = "160.16"
$ws.Range("E26").Value = "  +1.07%  "
$ws.Range("D27").Value = "20.61"
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").Value = "2.415"
$ws.Range("E28").Value = "  -2.59%  "
$ws.Range("D29").Value = "126.34"
$ws.Range("E29").Value = "  +1.04%  "
$ws.Range("D30").Value = "0.1057"
$ws.Range("E30").Value = "  -0.54%  "
$ws.Range("D31").Value = "1.033"
$ws.Range("E31").Value = "  -0.71%  "
$ws.Range("D32").Value = "5.845"
$ws.Range("E32").Value = "  -0.71%  "
$ws.Range("D33").Value = "3.603"
$ws.Range("E33").Value = "  +0.31%  "
$ws.Range("D34").Value = "0.02428"
$ws.Range("E34").Value = "  -0.36%  "
$ws.Range("D35").Value = "0.06473"
$ws.Range("E35").Value = "  -0.96%  "
$ws.Range("D36").Value = "9.075"
$ws.Range("E36").Value = "  -4.70%  "
$ws.Range("D37").Value = "0.2169"
$ws.Range("E37").Value = "  -0.70%  "
$ws.Range("D38").Value = "1.249"
$ws.Range("E38").Value = "  +1.79%  "
$ws.Range("D39").Value = "0.6435"
$ws.Range("E39").Value = "  -1.08%  "
$ws.Range("D40").Value = "1.182"
$ws.Range("E40").Value = "  -1.90%  "
$ws.Range("D41").Value = "4.944"
$ws.Range("E41").Value = "  -1.05%  "
$ws.Range("D42").Value = "11.07"
$ws.Range("E42").Value = "  -2.27%  "
$ws.Range("D43").Value = "0.6006"
$ws.Range("D44").Value = "12.89"
$ws.Range("E44").Value = "  -0.83%  "
$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "1.276"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "3.668"
$ws.Range("E46").Value = "  -0.16%  "
$ws.Range("D47").Value = "2.005"
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("D48").Value = "1.205"
$ws.Range("E48").Value = "  -0.95%  "
$ws.Range("D49").Value = "121.13"
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("D50").Value = "0.06869"
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("D51").Value = "77.08"
$ws.Range("E51").Value = "  -0.99%  "
